$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input 1 - Normal")
$ws.Activate()

# Fill in the running-number column (G) for the "TASK" block (rows 3-33),
# numbering starts at 1 on row 3 and increases by one per row.
for ($r = 3; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = ($r - 2)
    # The source cells carry no explicit style (General/Normal) even when
    # the row itself is custom-formatted, so strip any inherited style.
    $cell.Style = "Normal"
}

# Update the view: scroll position and the active selection.
$ws.Range("A76").Select()
$ws.Range("E30").Select()
